$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the existing date-number-format style (currently numFmtId 14, used by Q1)
# to a custom format "d/m/yy;@" so the same style can be reused by the new F2 cell.
$ws.Range("Q1").NumberFormat = "d/m/yy;@"

# Add the new data row (row 2)
$ws.Range("A2").Value = 123
$ws.Range("B2").Value = 456
$ws.Range("C2").Value = "encapsulant"
$ws.Range("D2").Value = "DH"
$ws.Range("E2").Value = 200

# F2 holds a date (2025-04-02), formatted the same way as Q1
$ws.Range("F2").Value = 45749
$ws.Range("F2").NumberFormat = "d/m/yy;@"

$ws.Range("Q2").Value = 24

# Widen column F to fit the new date values/header
$ws.Columns.Item(6).ColumnWidth = 10.28

# Update the active selection / scroll position shown in the sheet view
$null = $ws.Range("Q4").Select()
